# ---------------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx : table style swap + theme colour re-palette
# ---------------------------------------------------------------------------
# 1. Slide 16's summary table switches from the deck's custom "Table_0"
#    style ({D0AE3B76-4DBC-4395-875A-A3B0B0DB1FBA}) to the built-in table
#    style {5EF94A35-1368-4C70-A22B-DF403873E76B}. Table styles are not a
#    plain read/write property on the Table object (PowerPoint throws
#    "Table styles cannot be assigned through a property" if you try), so
#    this has to go through Table.ApplyStyle(...).
# 2. The deck's main theme (ppt/theme/theme1.xml, "Integral") is re-coloured
#    to match the Office default theme palette that used to live in the
#    (otherwise unused-by-slides) notes-master theme. We drive this through
#    the 12-slot ThemeColorScheme exposed on a slide, which writes straight
#    back into the shared theme part used by every slide's master.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------

$tableSlide = $null
$tableShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tableSlide = $slide
            $tableShape = $shape
        }
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{5EF94A35-1368-4C70-A22B-DF403873E76B}")
}

# --- 2. Theme colour scheme ------------------------------------------------

function Set-ThemeColorHex($themeColorScheme, $slotIndex, $hexRgb) {
    $r = [Convert]::ToInt32($hexRgb.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hexRgb.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hexRgb.Substring(4, 2), 16)
    $themeColorScheme.Colors($slotIndex).RGB = $r + ($g * 256) + ($b * 65536)
}

# Slot order (matches MsoThemeColorSchemeIndex / the 12 clrScheme children):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officePalette = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
foreach ($slot in $officePalette.Keys) {
    Set-ThemeColorHex $themeColorScheme $slot $officePalette[$slot]
}
